$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Training Dashboard")

# Helper cell far outside the used range (A1:K34) used to stage a pure-text
# value ("04-Nov-2025") so that PasteSpecial(xlPasteValues) can drop the
# literal text into column I without Excel's date auto-recognition kicking
# in and without disturbing the destination cells' existing number format/style.
$helper = $ws.Cells.Item(100, 100)
$helper.Value2 = "'04-Nov-2025"
$helper.Copy()

for ($r = 3; $r -le 34; $r++) {
    $hCell = $ws.Cells.Item($r, 8)   # column H - PERIOD TO EXPIRE
    $iCell = $ws.Cells.Item($r, 9)   # column I - LAST UPDATE

    $hCell.Value2 = $hCell.Value2 - 1
    $iCell.PasteSpecial(-4163)       # xlPasteValues: value only, keep style
}

$helper.Clear()
$excel.CutCopyMode = $false
